$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2 is formatted as Text (numFmtId 49), so set it explicitly so Excel
# stores the value as a shared string rather than a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "25/10/2022"
